# The sheet gained one new data row. A new row is inserted at sheet row 631
# (pushing the former rows 631-685 down to 632-686 unchanged) and the newly
# inserted row 631 is populated with its own record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(631).Insert()

$ws.Cells.Item(631, 1).Value = 3
$ws.Cells.Item(631, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(631, 3).Value = "Coquimbo"
$ws.Cells.Item(631, 4).Value = 45132
$ws.Cells.Item(631, 5).Value = 5
$ws.Cells.Item(631, 6).Value = 100112017
$ws.Cells.Item(631, 7).Value = "Apio"
$ws.Cells.Item(631, 8).Value = "Americana (o)"
$ws.Cells.Item(631, 9).Value = "Primera"
$ws.Cells.Item(631, 10).Value = 120
$ws.Cells.Item(631, 11).Value = 8000
$ws.Cells.Item(631, 12).Value = 8000
$ws.Cells.Item(631, 13).Value = 8000
$ws.Cells.Item(631, 14).Value = "$/docena de matas"
$ws.Cells.Item(631, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(631, 16).Value = 1333
$ws.Cells.Item(631, 17).Value = 6
$ws.Cells.Item(631, 18).Value = "Hortaliza"
